# Adds the "Off board components" section (rows 28-36) to the V2_1-BOM
# worksheet, including a real hyperlink on D30, then restores the
# selection to C31 as in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header ------------------------------------------------------
$ws.Range("A28").Value = "Off board components"
$ws.Range("A28").Font.Bold = $true

# --- Row 29: RGB TFT display ---------------------------------------------
$ws.Range("A29").Value = "RGB TFT display"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = "Waveshare"
$ws.Range("D29").Value = "https://www.amazon.nl/dp/B08H24H7KX?ref=ppx_yo2ov_dt_b_fed_asin_title&th=1"
$ws.Range("E29").Value = "2.4""SPI 240x320 with ILI9341 driver"

# --- Row 30: Heater element (with real hyperlink on D30) -----------------
$ws.Range("A30").Value = "Heater element"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "LJXH 170X62X5.5mm PTC"
$ws.Range("D30").Value = "https://nl.aliexpress.com/item/32854686343.html?spm=a2g0o.order_list.order_list_main.5.63bb79d2zTMKcl&gatewayAdapt=glo2nld"
$ws.Hyperlinks.Add($ws.Range("D30"), "https://nl.aliexpress.com/item/32854686343.html?spm=a2g0o.order_list.order_list_main.5.63bb79d2zTMKcl&gatewayAdapt=glo2nld")
$ws.Range("E30").Value = "220V 250 degreeC 300W"

# --- Row 31: Solid State Relais -------------------------------------------
$ws.Range("A31").Value = "Solid State Relais"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = "SSR-40 DA"
$ws.Range("D31").Value = "https://www.amazon.nl/Relais-Eenfase-SSR-40DA-24-380V-Relaismodule/dp/B07HB3JXGG/ref=sr_1_1_sspa?__mk_nl_NL=%C3%85M%C3%85%C5%BD%C3%95%C3%91&crid=3C2JFANRODBPX&dib=eyJ2IjoiMSJ9.i-azeWWndKw1g8EuvXNmWCcApIIj3aQPny46-c4GTNzZUi7q6biab31AZMgssa2KXRzzl3V65jFQEDaJ8rdbjOZRAk3ihFhtlF4XLW1vY3ffW90ranGTFNSh6XcczIkmfCpwq0Hq4zycbadHxDAY80BR15hFA_I7BJSa40qii5FFbffI5Lvk_LytXBJGWNrtjQWenkOyilNJbfVYN_lTLh1e3XhKjcAzgSCLE6p-fIGn_Ejo0QRHPhTksMGPTkHLU0eXE6leXT_MeoJ556D4DUW_wtJSWrYa9H9PF37x2LE.soiOgJh5tEPfqYilWH7Tit3mIVLvl9i1c8s6EotjzCI&dib_tag=se&keywords=ssr40+DA&qid=1742133138&sprefix=ssr40+da%2Caps%2C81&sr=8-1-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9hdGY&psc=1"

# --- Row 32: K-Type sensor --------------------------------------------------
$ws.Range("A32").Value = "K-Type sensor"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "Must have lug!"
$ws.Range("D32").Value = "https://www.amazon.nl/dp/B0C9F187RQ?ref=ppx_yo2ov_dt_b_fed_asin_title&th=1"

# --- Row 33: Insulation blanket --------------------------------------------
$ws.Range("A33").Value = "Insulation blanket"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "FLASLD aluminium heatshield blanket"
$ws.Range("D33").Value = "https://www.amazon.nl/dp/B0B19RNB3T?ref_=ppx_hzod_title_dt_b_fed_asin_title_0_0&th=1"
$ws.Range("E33").Value = "The one I ordered is no longer available, but there are many others"

# --- Row 34: Standoff -------------------------------------------------------
$ws.Range("A34").Value = "Standoff"
$ws.Range("B34").Value = 2
$ws.Range("C34").Value = "8mm M3"

# --- Row 35: Knob ------------------------------------------------------------
$ws.Range("A35").Value = "Knob"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "For the rotary encoder"

# --- Row 36: DC chassis part --------------------------------------------------
$ws.Range("A36").Value = "DC chassis part"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = "For the fan(s)"

# --- Column width tweaks (widened to fit the new, longer text) ---------------
$ws.Columns.Item(3).ColumnWidth = 35.140625
$ws.Columns.Item(4).ColumnWidth = 63.42578125

# --- Restore the view / selection state --------------------------------------
$excel.Goto($ws.Range("A7"), $true)
$ws.Range("C31").Select()
